# Append two new filtered-feed rows (rows 90 and 91) to the "Filtered Feeds"
# sheet, mirroring the existing link/keywords/title layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        Row     = 90
        Url     = "https://www.genomeweb.com/cancer/natera-submits-signatera-cdx-bladder-cancer-fda-premarket-approval"
        Keyword = "CDx"
        Title   = "Natera Submits Signatera CDx for Bladder Cancer to FDA for Premarket Approval"
    },
    @{
        Row     = 91
        Url     = "https://www.360dx.com/cancer/natera-submits-signatera-cdx-bladder-cancer-fda-premarket-approval"
        Keyword = "CDx"
        Title   = "Natera Submits Signatera CDx for Bladder Cancer to FDA for Premarket Approval"
    }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)

    $aCell.Value = $item.Url
    $bCell.Value = $item.Keyword
    $cCell.Value = $item.Title

    # Turn column A into a real hyperlink, then reapply the built-in
    # "Hyperlink" cell style so it matches the other link cells.
    $ws.Hyperlinks.Add($aCell, $item.Url) | Out-Null
    $aCell.Style = "Hyperlink"
}

Write-Output "Added rows 90-91 with Natera Signatera CDx bladder cancer feed entries"
